# Fixed the shooting input bug
# Applies the textual edits described by the commit diff using Word COM
# interop (Find & Replace plus direct Range manipulation for the larger
# paragraph rewrite).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "...but I also didn't make some packages I used..."
#    -> "...but I also didn't code some of the packages I used..."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "but I also didn’t make some packages",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "but I also didn’t code some of the packages", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) "...such as Audio and UI." -> "...such as Audio and UI, but these
#    are made by me."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "such as Audio and UI.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "such as Audio and UI, but these are made by me.", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) "...were all imported from aforementioned sources." -> "...were all
#    imported from the aforementioned sources."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "were all imported from aforementioned sources.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "were all imported from the aforementioned sources.", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) Replace the short "First roadblock: Multiple PlayerInputs" heading
#    paragraph text with the full war-story paragraph about the bug.
# ---------------------------------------------------------------------
$newText = "The firing was working just fine, but then I set up the main menu (which took 5 minutes since I already made a template) and it was not working anymore. The player was moving according to the input, but it was not shooting.  Although I don’t know why it is like this, the likely culprit is the PlayerInput component I have on my Canvas (which was necessary for my menu) which meant I had two PlayerInput components in my scene at once which broke something (the Input System package is still relatively new so it’s still a bit funky to use at times). The movement input was using messages and worked fine, while the shooting was using with C# events and didn’t work. I couldn’t use messages for the shooting since I need to detect whether the input was pressed or released so I switched the PlayerInput to use C# events and made sure to change the movement input accordingly. After that, the shooting and movement were both working fine, at the expense of concise code since it probably tripled the amount of line needed for this."

$d.Content.Find.Execute(
    "First roadblock: Multiple PlayerInputs",
    $true, $false, $false, $false, $false, $true, 1, $false,
    $newText, 2) | Out-Null
